$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($r in 2..6) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    $cell.Value = 45175
}
